# This edit inserts a new data row into the worksheet at row 60 (pushing the
# existing row 60 and all rows below it down by one), and populates the new
# row with its own data. This matches the commit "Fruta / hortaliza, semanal"
# which adds a new weekly price observation to the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60; everything from old row 60 onward
# shifts down to row+1 (old row 60 -> new row 61, ..., old row 139 -> new row 140).
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with its values.
# Columns that are identical across all rows for this market/category keep the
# same constant values (A, B, C, E, F, G, H, I, N, Q, R); D, J, K, L, M, O, P
# carry the new observation's data.
$ws.Cells.Item(60, 1).Value2 = 3
$ws.Cells.Item(60, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(60, 3).Value2 = "Coquimbo"
$ws.Cells.Item(60, 4).Value2 = 44601
$ws.Cells.Item(60, 5).Value2 = 5
$ws.Cells.Item(60, 6).Value2 = 100112030
$ws.Cells.Item(60, 7).Value2 = "Poroto granado"
$ws.Cells.Item(60, 8).Value2 = "Sin especificar"
$ws.Cells.Item(60, 9).Value2 = "Primera"
$ws.Cells.Item(60, 10).Value2 = 65
$ws.Cells.Item(60, 11).Value2 = 23000
$ws.Cells.Item(60, 12).Value2 = 24000
$ws.Cells.Item(60, 13).Value2 = 23538
$ws.Cells.Item(60, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(60, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(60, 16).Value2 = 942
$ws.Cells.Item(60, 17).Value2 = 25
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"
